# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# header style already used by column H, and fill in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto I1/J1, then set their text so the
# new header cells share H1's formatting (bold, bordered, centered/top).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2-7.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 6
